# Updates the cryptos list: refreshed prices / 1h volume percentages, and
# re-ranks a couple of coin pairs whose order flipped (rows 12/13 and 24/25).
#
# Price strings in column D (e.g. "25.920.48", "0.0636") must stay text --
# Excel's normal Value-assignment auto-detects plain numeric-looking strings
# and coerces them to numbers, which would change both the stored type and
# the display. Force text via NumberFormat "@" before assigning, then reset
# the style back to "Normal" so no stray formatting is left on the cell.
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "25.920.48"
$ws.Range("E2").Value = "  -1.32%  "

Set-TextValue $ws.Range("D3") "1.636.58"
$ws.Range("E3").Value = "  -0.67%  "

$ws.Range("E4").Value = "  +0.08%  "

Set-TextValue $ws.Range("D5") "215.40"
$ws.Range("E5").Value = "  -0.83%  "

$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  -0.94%  "

Set-TextValue $ws.Range("D9") "0.0636"
$ws.Range("E9").Value = "  -0.27%  "

$ws.Range("E10").Value = "  -1.61%  "

Set-TextValue $ws.Range("D11") "0.0794"
$ws.Range("E11").Value = "  +0.01%  "

# Row 12 / Row 13 swap (Polkadot now ranks above WrappedliquidstakedEther2.0)
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D12") "4.28"
$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D13") "1.863.17"
$ws.Range("E13").Value = "  -0.68%  "

Set-TextValue $ws.Range("D14") "1.627.92"
$ws.Range("E14").Value = "  -1.88%  "

Set-TextValue $ws.Range("D15") "0.543"
$ws.Range("E15").Value = "  -0.48%  "

Set-TextValue $ws.Range("D16") "0.0₃0762"
$ws.Range("E16").Value = "  -0.36%  "

Set-TextValue $ws.Range("D17") "62.82"
$ws.Range("E17").Value = "  -0.83%  "

Set-TextValue $ws.Range("D18") "25.929.80"

$ws.Range("E19").Value = "  +0.12%  "

Set-TextValue $ws.Range("D20") "192.73"
$ws.Range("E20").Value = "  -1.53%  "

Set-TextValue $ws.Range("D21") "4.35"
$ws.Range("E21").Value = "  -2.03%  "

$ws.Range("E23").Value = "  -0.83%  "

# Row 24 / Row 25 swap (Toncoin now ranks above Stellar)
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D24") "1.80"
$ws.Range("E24").Value = "  +0.70%  "

$ws.Range("B25").Value = "Stellar"
$ws.Range("C25").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D25") "0.131"
$ws.Range("E25").Value = "  +4.58%  "

Set-TextValue $ws.Range("D26") "143.87"
$ws.Range("E26").Value = "  +0.28%  "

$ws.Range("E27").Value = "  -0.03%  "

Set-TextValue $ws.Range("D28") "6.89"
$ws.Range("E28").Value = "  -0.80%  "

Set-TextValue $ws.Range("D29") "15.54"
$ws.Range("E29").Value = "  -0.70%  "

$ws.Range("E30").Value = "  -0.62%  "

$ws.Range("E31").Value = "  -0.49%  "

$ws.Range("E32").Value = "  -2.39%  "

Set-TextValue $ws.Range("D33") "3.25"
$ws.Range("E33").Value = "  -0.24%  "

$ws.Range("E34").Value = "  -4.43%  "

$ws.Range("E35").Value = "  +1.42%  "

$ws.Range("E36").Value = "  -1.40%  "

Set-TextValue $ws.Range("D37") "1.132.89"
$ws.Range("E37").Value = "  -0.26%  "

Set-TextValue $ws.Range("D38") "0.543"

$ws.Range("E39").Value = "  -1.89%  "

$ws.Range("E40").Value = "  -0.96%  "

$ws.Range("E41").Value = "  -0.42%  "

Set-TextValue $ws.Range("D42") "99.48"

Set-TextValue $ws.Range("D43") "0.793"
$ws.Range("E43").Value = "  -0.86%  "

Set-TextValue $ws.Range("D44") "1.772.77"
$ws.Range("E44").Value = "  -0.67%  "

Set-TextValue $ws.Range("D45") "0.0₆0115"
$ws.Range("E45").Value = "  +2.57%  "

Set-TextValue $ws.Range("D46") "56.59"
$ws.Range("E46").Value = "  -0.99%  "

$ws.Range("E47").Value = "  +2.25%  "

$ws.Range("E48").Value = "  -0.45%  "

Set-TextValue $ws.Range("D49") "7.71"
$ws.Range("E49").Value = "  +0.25%  "

$ws.Range("E50").Value = "  -0.92%  "

Set-TextValue $ws.Range("D51") "0.0960"
$ws.Range("E51").Value = "  -1.27%  "
